$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$r = $ws.Range("J2")
$r.Borders.LineStyle = -4142
